# REPORTGEN-1086: fix part of issue
#
# Updates the "ISO-5055 Full Detailed Report - OMG Technical Debt Edition"
# template so that:
#   - the Summary sheet's findings-summary caption / quality-standard label /
#     technical-debt-evolution RepGen tag are refreshed for the OMG -> ISO-5055
#     wording, and
#   - each "<Characteristic> Details" sheet's RepGen table tag for the rule
#     violations list gets an extra OMG=true parameter.
#
# NOTE: shared-string indices are NOT hard-coded anywhere below. The engine
# rebuilds xl/sharedStrings.xml (and every cell's numeric <v> index) from the
# live cell text on save, so we only need to set the handful of cells whose
# *text* actually changes, in the same order the strings should appear in the
# rebuilt table.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")

# "Findings summary for CAST under OMG-ASCQM Standards" -> "... ISO-5055 Standards"
$wsSummary.Range("B12").Value = "Findings summary for CAST under ISO-5055 Standards"

# RepGen:TABLE;OMG_TECHNICAL_DEBT_EVOLUTION;ID=ISO-5055-Index,HEADER=NO
#   -> adds MORE=true
$wsSummary.Range("B14").Value = "RepGen:TABLE;OMG_TECHNICAL_DEBT_EVOLUTION;ID=ISO-5055-Index,MORE=true,HEADER=NO"

# "Quality Standard" -> "ISO-5055"
$wsSummary.Range("B13").Value = "ISO-5055"

# Each "<Metric> Details" sheet's A2 RepGen table tag gains ",OMG=true"
$wsSecurityDetails = $wb.Worksheets.Item("ISO-5055-Security Details")
$wsSecurityDetails.Range("A2").Value = "RepGen:TABLE;LIST_RULES_VIOLATIONS_BOOKMARKS_TABLE;METRICS=ISO-5055-Security,COUNT=-1,HEADER=NO,OMG=true"

$wsReliabilityDetails = $wb.Worksheets.Item("ISO-5055-Reliability Details")
$wsReliabilityDetails.Range("A2").Value = "RepGen:TABLE;LIST_RULES_VIOLATIONS_BOOKMARKS_TABLE;METRICS=ISO-5055-Reliability,COUNT=-1,HEADER=NO,OMG=true"

$wsPerfEffDetails = $wb.Worksheets.Item("ISO-5055-Perf-Eff Details")
$wsPerfEffDetails.Range("A2").Value = "RepGen:TABLE;LIST_RULES_VIOLATIONS_BOOKMARKS_TABLE;METRICS=ISO-5055-Performance-Efficiency,COUNT=-1,HEADER=NO,OMG=true"

$wsMaintainabilityDetails = $wb.Worksheets.Item("ISO-505-Maintainability Details")
$wsMaintainabilityDetails.Range("A2").Value = "RepGen:TABLE;LIST_RULES_VIOLATIONS_BOOKMARKS_TABLE;METRICS=ISO-5055-Maintainability,COUNT=-1,HEADER=NO,OMG=true"

# --- Selection bookkeeping: the edited detail sheets now show A2 selected
#     (they used to show A3) ---

$wsSecurityDetails.Activate() | Out-Null
$wsSecurityDetails.Range("A2").Select() | Out-Null

$wsPerfEffDetails.Activate() | Out-Null
$wsPerfEffDetails.Range("A2").Select() | Out-Null

$wsMaintainabilityDetails.Activate() | Out-Null
$wsMaintainabilityDetails.Range("A2").Select() | Out-Null

# --- Restore Summary as the active sheet/tab, with its updated selection ---

$wsSummary.Activate() | Out-Null
$wsSummary.Range("B28").Select() | Out-Null
